$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = 0.02394092632924144
$ws.Range("B2").Value = 0
$ws.Range("C2").Value = -0
$ws.Range("D2").Value = -0.04737869069218277
$ws.Range("E2").Value = -0.0687263685828748
$ws.Range("F2").Value = -0
$ws.Range("J2").Value = 0.07200801432171179
$ws.Range("K2").Value = -0
$ws.Range("L2").Value = -0
$ws.Range("M2").Value = 0.01396478623361789
$ws.Range("N2").Value = -0.06907899306029104
$ws.Range("O2").Value = -0
$ws.Range("Q2").Value = -0
$ws.Range("R2").Value = 0
$ws.Range("S2").Value = -0.02592058921947108
$ws.Range("T2").Value = -0
$ws.Range("U2").Value = -0
$ws.Range("V2").Value = 0.09340842381038526
$ws.Range("W2").Value = 0.03130339279830054
$ws.Range("AA2").Value = 0
$ws.Range("AB2").Value = 0.002851801335650072
$ws.Range("AC2").Value = -0
$ws.Range("AD2").Value = -0
$ws.Range("AE2").Value = -0.0344727708704435
$ws.Range("AF2").Value = 0.01411038664133246
$ws.Range("AG2").Value = 0
$ws.Range("AI2").Value = 0
$ws.Range("AK2").Value = 0.02309732901816538
$ws.Range("AL2").Value = 0
$ws.Range("AM2").Value = 0
$ws.Range("AN2").Value = 0.01486262558060095
$ws.Range("AO2").Value = -0.06283453873230342
$ws.Range("AP2").Value = 0
$ws.Range("AQ2").Value = -0
$ws.Range("AR2").Value = 0
$ws.Range("AT2").Value = 0.01438911099255462
$ws.Range("AU2").Value = 0
$ws.Range("AV2").Value = 0
$ws.Range("AW2").Value = -0.02276828455198836
$ws.Range("AX2").Value = -0.01183823534383922
$ws.Range("BC2").Value = 0.03155251487596699
$ws.Range("BD2").Value = 0
$ws.Range("BE2").Value = -0
$ws.Range("BF2").Value = 0.002581439954956081
$ws.Range("BG2").Value = -0.02235394201285093
$ws.Range("BI2").Value = 0
$ws.Range("BK2").Value = -0
$ws.Range("BL2").Value = 0.03087649112700357
$ws.Range("BM2").Value = 0
$ws.Range("BN2").Value = -0
$ws.Range("BO2").Value = 0.01162501093338994
$ws.Range("BP2").Value = 0.008739102550135084
$ws.Range("BQ2").Value = -0
$ws.Range("BS2").Value = -0
$ws.Range("BT2").Value = -0
$ws.Range("BU2").Value = 0.01026856340313644
$ws.Range("BV2").Value = -0
$ws.Range("BW2").Value = -0
$ws.Range("BX2").Value = 0.0007358185974277789
$ws.Range("BY2").Value = -0.04082188383782141
$ws.Range("BZ2").Value = -0
$ws.Range("CD2").Value = -0.01118785406560184
$ws.Range("CE2").Value = 0
$ws.Range("CF2").Value = 0
$ws.Range("CG2").Value = -0.007797836553935383
$ws.Range("CH2").Value = 0.02072217619577005
$ws.Range("CK2").Value = -0
$ws.Range("CL2").Value = -0
$ws.Range("CM2").Value = -0.0037524436061281
$ws.Range("CN2").Value = 0
$ws.Range("CO2").Value = -0
$ws.Range("CP2").Value = -0.0250987332193468
$ws.Range("CQ2").Value = -0.04909057745741415
$ws.Range("CR2").Value = 0
$ws.Range("CS2").Value = -0
$ws.Range("CV2").Value = -0.006133870860311784
$ws.Range("CW2").Value = -0
$ws.Range("CX2").Value = -0
$ws.Range("CY2").Value = 0.01491518537000321
$ws.Range("CZ2").Value = 0.03547509364258496
$ws.Range("DC2").Value = 0
$ws.Range("DD2").Value = -0
$ws.Range("DE2").Value = 0.002217122834552273
$ws.Range("DF2").Value = 0
$ws.Range("DG2").Value = 0
$ws.Range("DH2").Value = 0.06416965779138531
$ws.Range("DI2").Value = -0.05585497434061494
$ws.Range("DJ2").Value = -0
$ws.Range("DN2").Value = 0.002884060647310469
$ws.Range("DO2").Value = -0
$ws.Range("DP2").Value = -0
$ws.Range("DQ2").Value = 0.01895263501235869
$ws.Range("DR2").Value = 0.1034217657051097
$ws.Range("DT2").Value = -0
$ws.Range("DU2").Value = -0
$ws.Range("DW2").Value = 0.00729196034322644
$ws.Range("DX2").Value = 0
$ws.Range("DY2").Value = 0
$ws.Range("DZ2").Value = -0.008249161542692014
$ws.Range("EA2").Value = 0.03280690868587788
$ws.Range("EB2").Value = -0
$ws.Range("EC2").Value = -0
$ws.Range("ED2").Value = 0
$ws.Range("EF2").Value = -0.01677180019557874
$ws.Range("EG2").Value = 0
$ws.Range("EH2").Value = 0
$ws.Range("EI2").Value = 0.03565718356415532
$ws.Range("EJ2").Value = 0.08419176871232947
$ws.Range("EK2").Value = 0
$ws.Range("EL2").Value = -0
$ws.Range("EO2").Value = -0.02802871556604173
$ws.Range("EP2").Value = -0
$ws.Range("EQ2").Value = 0
$ws.Range("ER2").Value = 0.006275292382618154
$ws.Range("ES2").Value = 0.05217510566426688
$ws.Range("EV2").Value = 0
$ws.Range("EW2").Value = -0
$ws.Range("EX2").Value = -0.008837534930016837
$ws.Range("EY2").Value = 0
$ws.Range("EZ2").Value = 0
$ws.Range("FA2").Value = 0.007528061814176132
$ws.Range("FB2").Value = 0.03430455635138342
$ws.Range("FE2").Value = -0
$ws.Range("FF2").Value = -0
$ws.Range("FG2").Value = 0.01735606272529178
$ws.Range("FH2").Value = 0
$ws.Range("FI2").Value = -0
$ws.Range("FJ2").Value = 0.05155398105349183
$ws.Range("FK2").Value = 0.06938520987199787
$ws.Range("FM2").Value = -0
$ws.Range("FN2").Value = -0
$ws.Range("FP2").Value = -0.0195193089742912
$ws.Range("FQ2").Value = 0
$ws.Range("FR2").Value = 0
$ws.Range("FS2").Value = -0.03508992439640592
$ws.Range("FT2").Value = 0.1510505264966768
$ws.Range("FU2").Value = -0
$ws.Range("FX2").Value = -0
$ws.Range("FY2").Value = -0.003483313476346797
$ws.Range("FZ2").Value = 0
$ws.Range("GA2").Value = -0
$ws.Range("GB2").Value = 0.03446025714747364
$ws.Range("GC2").Value = 0
$ws.Range("GD2").Value = 0
$ws.Range("GE2").Value = -0
$ws.Range("GG2").Value = -0
$ws.Range("A3").Value = -0
$ws.Range("C3").Value = -0
$ws.Range("D3").Value = -0.0630659001931053
$ws.Range("E3").Value = -0.01300764756401055
$ws.Range("F3").Value = 0.5240954286294019
$ws.Range("G3").Value = -0
$ws.Range("I3").Value = 0
$ws.Range("J3").Value = -0
$ws.Range("K3").Value = 0
$ws.Range("L3").Value = -0
$ws.Range("M3").Value = -0.05405252220731481
$ws.Range("N3").Value = -0.05525720831395486
$ws.Range("O3").Value = 0.3895578008224662
$ws.Range("Q3").Value = 0
$ws.Range("R3").Value = 0
$ws.Range("S3").Value = 0
$ws.Range("T3").Value = -0
$ws.Range("U3").Value = 0
$ws.Range("V3").Value = 0.01643674394651077
$ws.Range("W3").Value = 0.06939873613788489
$ws.Range("X3").Value = 0.04519975493257287
$ws.Range("Y3").Value = -0
$ws.Range("Z3").Value = 0
$ws.Range("AA3").Value = -0
$ws.Range("AB3").Value = -0
$ws.Range("AE3").Value = -0.0539946543198986
$ws.Range("AF3").Value = -0.000731475144162448
$ws.Range("AG3").Value = 0.03855298895738098
$ws.Range("AJ3").Value = 0
$ws.Range("AK3").Value = -0
$ws.Range("AN3").Value = -0.03828147137014404
$ws.Range("AO3").Value = -0.08486099760458013
$ws.Range("AP3").Value = 0.05670348389126057
$ws.Range("AQ3").Value = -0
$ws.Range("AR3").Value = 0
$ws.Range("AS3").Value = -0
$ws.Range("AT3").Value = -0
$ws.Range("AW3").Value = -0.0313681254931311
$ws.Range("AX3").Value = 0.04310579127287367
$ws.Range("AY3").Value = 0.3325697888168643
$ws.Range("AZ3").Value = -0
$ws.Range("BB3").Value = -0
$ws.Range("BC3").Value = 0
$ws.Range("BF3").Value = -0.01209553597845399
$ws.Range("BG3").Value = -0.009287378463346888
$ws.Range("BH3").Value = 0.07097712646542201
$ws.Range("BJ3").Value = 0
$ws.Range("BK3").Value = 0
$ws.Range("BL3").Value = -0
$ws.Range("BM3").Value = 0
$ws.Range("BN3").Value = 0
$ws.Range("BO3").Value = 0.02733952297264467
$ws.Range("BP3").Value = -0.00559636509108988
$ws.Range("BQ3").Value = 0.03514048739457296
$ws.Range("BR3").Value = -0
$ws.Range("BS3").Value = -0
$ws.Range("BT3").Value = -0
$ws.Range("BU3").Value = -0
$ws.Range("BV3").Value = -0
$ws.Range("BX3").Value = -0.0401298438780733
$ws.Range("BY3").Value = 0.07580772388942704
$ws.Range("BZ3").Value = 0.1560123355299239
$ws.Range("CA3").Value = -0
$ws.Range("CB3").Value = 0
$ws.Range("CC3").Value = -0
$ws.Range("CD3").Value = 0
$ws.Range("CF3").Value = 0
$ws.Range("CG3").Value = 0.02159828231322356
$ws.Range("CH3").Value = 0.00867571310732208
$ws.Range("CI3").Value = -0.1023849726231967
$ws.Range("CJ3").Value = 0
$ws.Range("CL3").Value = 0
$ws.Range("CM3").Value = -0
$ws.Range("CO3").Value = 0
$ws.Range("CP3").Value = -0.03771215862096453
$ws.Range("CQ3").Value = -0.09500487350482544
$ws.Range("CR3").Value = -0.009261628607069274
$ws.Range("CT3").Value = -0
$ws.Range("CU3").Value = 0
$ws.Range("CV3").Value = -0
$ws.Range("CY3").Value = 0.02644305643145452
$ws.Range("CZ3").Value = 0.01555079161548294
$ws.Range("DA3").Value = -0.09727711389827733
$ws.Range("DD3").Value = -0
$ws.Range("DE3").Value = 0
$ws.Range("DG3").Value = 0
$ws.Range("DH3").Value = 0.03258621407065847
$ws.Range("DI3").Value = -0.013928588507449
$ws.Range("DJ3").Value = -0.03128692613813622
$ws.Range("DK3").Value = 0
$ws.Range("DM3").Value = -0
$ws.Range("DN3").Value = -0
$ws.Range("DO3").Value = 0
$ws.Range("DQ3").Value = -0.04395409460582739
$ws.Range("DR3").Value = 0.04710913878632132
$ws.Range("DS3").Value = 0.01580127312378971
$ws.Range("DU3").Value = 0
$ws.Range("DV3").Value = 0
$ws.Range("DW3").Value = 0
$ws.Range("DY3").Value = 0
$ws.Range("DZ3").Value = -0.007316915824470367
$ws.Range("EA3").Value = 0.007557786809462135
$ws.Range("EB3").Value = -0.06475067601211205
$ws.Range("EE3").Value = -0
$ws.Range("EF3").Value = 0
$ws.Range("EG3").Value = -0
$ws.Range("EH3").Value = 0
$ws.Range("EI3").Value = 0.085015046155024
$ws.Range("EJ3").Value = 0.01319452892834481
$ws.Range("EK3").Value = -0.0900866339165938
$ws.Range("EN3").Value = 0
$ws.Range("EO3").Value = -0
$ws.Range("ER3").Value = 0.05886629268000991
$ws.Range("ES3").Value = 0.04944761632630492
$ws.Range("ET3").Value = -0.0718264947522245
$ws.Range("EU3").Value = -0
$ws.Range("EW3").Value = -0
$ws.Range("EX3").Value = 0
$ws.Range("EZ3").Value = 0
$ws.Range("FA3").Value = 0.0188312515658069
$ws.Range("FB3").Value = -0.004237364955702617
$ws.Range("FC3").Value = -0.09411699122155258
$ws.Range("FF3").Value = 0
$ws.Range("FG3").Value = -0
$ws.Range("FH3").Value = -0
$ws.Range("FJ3").Value = -0.005822987135418163
$ws.Range("FK3").Value = -0.03443632606135918
$ws.Range("FL3").Value = 0.02510756393004237
$ws.Range("FO3").Value = 0
$ws.Range("FP3").Value = 0
$ws.Range("FQ3").Value = 0
$ws.Range("FR3").Value = 0
$ws.Range("FS3").Value = -0.01893127849923146
$ws.Range("FT3").Value = -0.1561887957586328
$ws.Range("FU3").Value = -0.07997069028232594
$ws.Range("FV3").Value = 0
$ws.Range("FW3").Value = 0
$ws.Range("FX3").Value = -0
$ws.Range("FY3").Value = -0
$ws.Range("GB3").Value = 0.04949523247384263
$ws.Range("GC3").Value = -0
$ws.Range("GD3").Value = 0.04010624219646641
$ws.Range("GG3").Value = -0
